# Auto-generated edit applying the diff to Lamia_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3194.7856
$ws.Range("I33").Value = 192.9
$ws.Range("K33").Value = 192.9
$ws.Range("M33").Value = 36.09999999999999

$ws.Range("H64").Value = 6556.1113
$ws.Range("I64").Value = 2501
$ws.Range("J64").Value = 7063
$ws.Range("K64").Value = 2501
$ws.Range("L64").Value = 7063
$ws.Range("M64").Value = -2253
$ws.Range("N64").Value = -7559

$ws.Range("H67").Value = 6556.1113
$ws.Range("I67").Value = 2501
$ws.Range("J67").Value = 7063
$ws.Range("K67").Value = 2501
$ws.Range("L67").Value = 7063
$ws.Range("M67").Value = -1643
$ws.Range("N67").Value = -8779

$ws.Range("H70").Value = 6063665.5
$ws.Range("J70").Value = 7146039
$ws.Range("L70").Value = 21438117
$ws.Range("N70").Value = -21438657

$ws.Range("H73").Value = 6063665.5
$ws.Range("J73").Value = 7146039
$ws.Range("L73").Value = 21438117
$ws.Range("N73").Value = -21439989

$ws.Range("H74").Value = 9312.9375
$ws.Range("I74").Value = 10667.167
$ws.Range("J74").Value = 8500.4
$ws.Range("K74").Value = 10667.167
$ws.Range("L74").Value = 8500.4
$ws.Range("M74").Value = -9731.166999999999
$ws.Range("N74").Value = -10372.4

$ws.Range("H77").Value = 9312.9375
$ws.Range("I77").Value = 10667.167
$ws.Range("J77").Value = 8500.4
$ws.Range("K77").Value = 53335.835
$ws.Range("L77").Value = 42502
$ws.Range("M77").Value = -48655.835
$ws.Range("N77").Value = -51862

$ws.Range("H127").Value = 438.7143
$ws.Range("I127").Value = 438.7143
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1316.1429
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 3643.8571
$ws.Range("N127").ClearContents()

$ws.Range("H131").Value = 7155.5
$ws.Range("I131").Value = 6522.25
$ws.Range("K131").Value = 19566.75
$ws.Range("M131").Value = -14526.75

$ws.Range("H137").Value = 3442.0889
$ws.Range("I137").Value = 2756.138
$ws.Range("K137").Value = 8268.414000000001
$ws.Range("M137").Value = -5718.414000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 55558052
$ws.Range("I45").Value = 71430824
$ws.Range("K45").Value = 71430824
$ws.Range("M45").Value = -71430447

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H125").Value = 99999.5
$ws.Range("J125").Value = 99999.5
$ws.Range("L125").Value = 99999.5
$ws.Range("N125").Value = -109839.5

$ws.Range("H132").Value = 6706.375
$ws.Range("I132").Value = 4924.407
$ws.Range("K132").Value = 14773.221
$ws.Range("M132").Value = -12243.221

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28715.25
$ws.Range("I31").Value = 4302.625
$ws.Range("J31").Value = 34140.277
$ws.Range("K31").Value = 4302.625
$ws.Range("L31").Value = 34140.277
$ws.Range("M31").Value = -4007.625
$ws.Range("N31").Value = -34730.277

$ws.Range("H33").Value = 1999
$ws.Range("I33").Value = 1999
$ws.Range("K33").Value = 1999
$ws.Range("M33").Value = -1620

$ws.Range("H34").Value = 28715.25
$ws.Range("I34").Value = 4302.625
$ws.Range("J34").Value = 34140.277
$ws.Range("K34").Value = 4302.625
$ws.Range("L34").Value = 34140.277
$ws.Range("M34").Value = -4100.625
$ws.Range("N34").Value = -34544.277

$ws.Range("H47").Value = 34052
$ws.Range("J47").Value = 42105
$ws.Range("L47").Value = 42105
$ws.Range("N47").Value = -43237

$ws.Range("H134").Value = 3182.5312
$ws.Range("I134").Value = 1726.2084
$ws.Range("J134").Value = 7551.5
$ws.Range("K134").Value = 5178.6252
$ws.Range("L134").Value = 22654.5
$ws.Range("M134").Value = -2643.6252
$ws.Range("N134").Value = -27724.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 64102.24
$ws.Range("I2").Value = 17.545454
$ws.Range("J2").Value = 114454.5
$ws.Range("K2").Value = 105.272724
$ws.Range("L2").Value = 686727
$ws.Range("M2").Value = 7.727276000000003
$ws.Range("N2").Value = -686953

$ws.Range("H9").Value = 1371562.5
$ws.Range("I9").Value = 3650333.2
$ws.Range("J9").Value = 4300
$ws.Range("K9").Value = 10950999.6
$ws.Range("L9").Value = 12900
$ws.Range("M9").Value = -10950775.6
$ws.Range("N9").Value = -13348

$ws.Range("H54").Value = 6730.6665
$ws.Range("J54").Value = 6471.5
$ws.Range("L54").Value = 19414.5
$ws.Range("N54").Value = -20532.5

$ws.Range("H107").Value = 3901.75
$ws.Range("J107").Value = 5068
$ws.Range("L107").Value = 15204
$ws.Range("N107").Value = -19044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8117
$ws.Range("I70").Value = 7489.5
$ws.Range("K70").Value = 7489.5
$ws.Range("M70").Value = -7219.5

$ws.Range("H73").Value = 8117
$ws.Range("I73").Value = 7489.5
$ws.Range("K73").Value = 7489.5
$ws.Range("M73").Value = -6553.5

$ws.Range("H132").Value = 5547.7144
$ws.Range("I132").Value = 4515.6665
$ws.Range("J132").Value = 7405.4
$ws.Range("K132").Value = 13546.9995
$ws.Range("L132").Value = 22216.2
$ws.Range("M132").Value = -11016.9995
$ws.Range("N132").Value = -27276.2

$ws.Range("H134").Value = 31475
$ws.Range("J134").Value = 31475
$ws.Range("L134").Value = 94425
$ws.Range("N134").Value = -99495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4149.4
$ws.Range("I61").Value = 4149.4
$ws.Range("K61").Value = 4149.4
$ws.Range("M61").Value = -3947.4

$ws.Range("H68").Value = 4532.4644
$ws.Range("I68").Value = 3477.52
$ws.Range("J68").Value = 13323.667
$ws.Range("K68").Value = 3477.52
$ws.Range("L68").Value = 13323.667
$ws.Range("M68").Value = -2728.52
$ws.Range("N68").Value = -14821.667

$ws.Range("H71").Value = 4532.4644
$ws.Range("I71").Value = 3477.52
$ws.Range("J71").Value = 13323.667
$ws.Range("K71").Value = 17387.6
$ws.Range("L71").Value = 66618.33499999999
$ws.Range("M71").Value = -13643.6
$ws.Range("N71").Value = -74106.33499999999

$ws.Range("H113").Value = 4149.4
$ws.Range("I113").Value = 4149.4
$ws.Range("K113").Value = 4149.4
$ws.Range("M113").Value = -1979.4
